$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-11 Wednesday", "2024-09-12 Thursday"),
    @("85×57=", "85×75="),
    @("92×86=", "79×29="),
    @("36×98=", "60×74="),
    @("83×95=", "75×59="),
    @("49×73=", "71×80="),
    @("71×40=", "71×92="),
    @("50×85=", "88×79="),
    @("17×48=", "12×73="),
    @("33×65=", "15×47="),
    @("26×94=", "63×49="),
    @("67×15=", "27×15="),
    @("34×13=", "37×80="),
    @("60×13=", "21×19="),
    @("68×89=", "79×79="),
    @("51×55=", "84×37="),
    @("69×69=", "57×68="),
    @("54×18=", "81×65="),
    @("91×79=", "50×70="),
    @("67×89=", "45×59="),
    @("36×46=", "65×13="),
    @("34×81=", "96×95="),
    @("27×37=", "78×53="),
    @("59×56=", "84×75="),
    @("93×99=", "20×47="),
    @("94×61=", "51×34=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
